# Update the worksheet date and the twenty-five division problems.
#
# Cell text is updated via Cell.Range.Text (which replaces just the
# content and preserves the trailing paragraph/cell marks along with the
# existing run formatting), addressed by explicit row/column, instead of
# a blanket Find/Replace. Several "old" values reappear as "new" values
# elsewhere in the table (e.g. 31÷7= and 47÷9=), so a global
# Find-and-ReplaceAll could cascade and corrupt unrelated cells; direct
# per-cell assignment avoids that ambiguity entirely.

$d = $word.ActiveDocument

# Title line with the date.
$d.Paragraphs.Item(1).Range.Text = "2023-11-10 Friday"

$t = $d.Tables.Item(1)

# Map of (row, column) -> new value, in document order, for the five
# data rows of the table (rows 1, 5, 9, 13, 17; the rows in between hold
# blank student-work cells).
$updates = @(
    @{ Row = 1;  Col = 1; Text = "31÷9=" },
    @{ Row = 1;  Col = 2; Text = "44÷5=" },
    @{ Row = 1;  Col = 3; Text = "40÷3=" },
    @{ Row = 1;  Col = 4; Text = "35÷4=" },
    @{ Row = 1;  Col = 5; Text = "13÷8=" },

    @{ Row = 5;  Col = 1; Text = "18÷4=" },
    @{ Row = 5;  Col = 2; Text = "47÷9=" },
    @{ Row = 5;  Col = 3; Text = "15÷8=" },
    @{ Row = 5;  Col = 4; Text = "17÷9=" },
    @{ Row = 5;  Col = 5; Text = "72÷4=" },

    @{ Row = 9;  Col = 1; Text = "31÷7=" },
    @{ Row = 9;  Col = 2; Text = "99÷9=" },
    @{ Row = 9;  Col = 3; Text = "94÷9=" },
    @{ Row = 9;  Col = 4; Text = "80÷6=" },
    @{ Row = 9;  Col = 5; Text = "64÷8=" },

    @{ Row = 13; Col = 1; Text = "65÷7=" },
    @{ Row = 13; Col = 2; Text = "19÷6=" },
    @{ Row = 13; Col = 3; Text = "47÷3=" },
    @{ Row = 13; Col = 4; Text = "48÷4=" },
    @{ Row = 13; Col = 5; Text = "25÷9=" },

    @{ Row = 17; Col = 1; Text = "20÷8=" },
    @{ Row = 17; Col = 2; Text = "15÷2=" },
    @{ Row = 17; Col = 3; Text = "44÷4=" },
    @{ Row = 17; Col = 4; Text = "12÷2=" },
    @{ Row = 17; Col = 5; Text = "33÷5=" }
)

foreach ($u in $updates) {
    $t.Cell($u.Row, $u.Col).Range.Text = $u.Text
}
